# Applies odds updates to Sheet1 (rows 5, 6 and 7) as described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7

# Row 6
$ws.Range("J6").Value = 2.5
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("X6").Value = 7
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 13
$ws.Range("AA6").Value = 17
$ws.Range("AH6").Value = 23

# Row 7
$ws.Range("G7").Value = 1.44
$ws.Range("H7").Value = 4.33
$ws.Range("I7").Value = 7.5
$ws.Range("J7").Value = 2
$ws.Range("L7").Value = 7
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 9.5
$ws.Range("Q7").Value = 1.98
$ws.Range("R7").Value = 1.88
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.75
$ws.Range("Z7").Value = 9
$ws.Range("AC7").Value = 9.5
$ws.Range("AD7").Value = 8
$ws.Range("AE7").Value = 21
$ws.Range("AF7").Value = 67
$ws.Range("AI7").Value = 21
$ws.Range("AJ7").Value = 81
$ws.Range("AL7").Value = 51
$ws.Range("AN7").Value = 3.25
$ws.Range("AO7").Value = 7
$ws.Range("AQ7").Value = 21
$ws.Range("AR7").Value = 51
$ws.Range("AT7").Value = 2.75
$ws.Range("AU7").Value = 9.5
$ws.Range("AW7").Value = 8.5
$ws.Range("AZ7").Value = 151
